$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row 15 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(15, 1).Value = "Afmelding nieuwsbrief"
$logs.Cells.Item(15, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(15, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Cells.Item(15, 4).Value = "Afmelding / Nieuwsbrief"
$logs.Cells.Item(15, 6).Value = "2025-06-20 14:30:44"
$logs.Cells.Item(15, 7).Value = "Nee"

# Extend the conditional formatting ranges on columns D and G so they
# cover the newly added row (D2:D14 -> D2:D15, G2:G14 -> G2:G15).
$rngD = $logs.Range("D2:D15")
$rngD.FormatConditions.Item(1).ModifyAppliesToRange($rngD)

$rngG = $logs.Range("G2:G15")
$rngG.FormatConditions.Item(1).ModifyAppliesToRange($rngG)

# --- Dashboard sheet: swap the first two category rows / update counts ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(2, 1).Value = "Afmelding / Nieuwsbrief"
$dash.Cells.Item(2, 2).Value = 4

$dash.Cells.Item(3, 1).Value = "Samenwerking / Partnerverzoek"
$dash.Cells.Item(3, 2).Value = 4
